$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Insert 5 new columns before column E (old E:I shifts to J:N)
$ws.Range("E1:I1").EntireColumn.Insert()

# ---- Row 8 header (first table) ----
$ws.Range("E8").Value = "فصل اول منتهی به 1399/03"
$ws.Range("F8").Value = "فصل دوم منتهی به 1399/06"
$ws.Range("G8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("H8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("I8").Value = "فصل اول منتهی به 1400/03"

# ---- Row 24 header (second table) ----
$ws.Range("E24").Value = "فصل اول منتهی به 1399/03"
$ws.Range("F24").Value = "فصل دوم منتهی به 1399/06"
$ws.Range("G24").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("H24").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("I24").Value = "فصل اول منتهی به 1400/03"

# ---- Row 10: هزینه حمل و نقل و انتقال ----
$ws.Range("E10").Value = 6013885
$ws.Range("F10").Value = 6160019
$ws.Range("G10").Value = 10783849
$ws.Range("H10").Value = 8399574
$ws.Range("I10").Value = 8033579

# ---- Row 11: zero row ----
$ws.Range("E11:I11").Value = 0

# ---- Row 12: zero row ----
$ws.Range("E12:I12").Value = 0

# ---- Row 13: zero row ----
$ws.Range("E13:I13").Value = 0

# ---- Row 14: هزینه تبلیغات ----
$ws.Range("E14").Value = 6203
$ws.Range("F14").Value = 2647
$ws.Range("G14").Value = 10872
$ws.Range("H14").Value = 15478
$ws.Range("I14").Value = 9159

# ---- Row 15: zero row ----
$ws.Range("E15:I15").Value = 0

# ---- Row 16: هزینه استهلاک ----
$ws.Range("E16").Value = 4847
$ws.Range("F16").Value = 4846
$ws.Range("G16").Value = 4847
$ws.Range("H16").Value = 155073
$ws.Range("I16").Value = 89032

# ---- Row 17: هزینه حقوق و دستمزد ----
$ws.Range("E17").Value = 87468
$ws.Range("F17").Value = 146276
$ws.Range("G17").Value = 103911
$ws.Range("H17").Value = 150168
$ws.Range("I17").Value = 200712

# ---- Row 18: zero row ----
$ws.Range("E18:I18").Value = 0

# ---- Row 19: سایر هزینه ها ----
$ws.Range("E19").Value = 439354
$ws.Range("F19").Value = 267972
$ws.Range("G19").Value = 614130
$ws.Range("H19").Value = 629700
$ws.Range("I19").Value = 519399

# ---- Row 20: جمع ----
$ws.Range("E20").Value = 6551757
$ws.Range("F20").Value = 6581760
$ws.Range("G20").Value = 11517609
$ws.Range("H20").Value = 9349993
$ws.Range("I20").Value = 8851881

# ---- Row 26: تعداد پرسنل غیر تولیدی شرکت ----
$ws.Range("E26").Value = 703
$ws.Range("F26").Value = 705
$ws.Range("G26").Value = 701
$ws.Range("H26").Value = 702
$ws.Range("I26").Value = 701

# ---- Row 27: تعداد پرسنل تولیدی شرکت ----
$ws.Range("E27").Value = 104
$ws.Range("F27").Value = 104
$ws.Range("G27").Value = 104
$ws.Range("H27").Value = 104
$ws.Range("I27").Value = 104
